# Datum_Conversion.xlsx update: remove three obsolete USGS/USACE station rows
# (82770, 76220, 76593) from Sheet1. Remaining rows shift up to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete whole rows 21, 17, 6 (highest row number first so earlier row
# numbers stay valid while we work).
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(6).Delete()

$ws.Range("C29").Select()
